$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Study Plan App - Burndown Chart")

# Update J9 (Pts Done Sprint 4.0 for task "9. Create Review page") from 0 to 2
$ws.Range("J9").Value = 2

# Update the selected cell/range in the sheet view
$ws.Range("J10").Select()
